# Insert a new data row for "Feria Lagunitas de Puerto Montt - Ciboulette" (Hortaliza)
# at row 237, shifting the existing rows 237:293 down to 238:294.
# The workbook's used range grows from A1:R293 to A1:R294.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 237 (pushes rows 237-293 down to 238-294)
$ws.Rows.Item(237).Insert()

# Populate the new row 237 with the new record's data
$row = 237
$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 44932
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100112039
$ws.Cells.Item($row, 7).Value = "Ciboulette"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 240
$ws.Cells.Item($row, 11).Value = 3500
$ws.Cells.Item($row, 12).Value = 3500
$ws.Cells.Item($row, 13).Value = 3500
$ws.Cells.Item($row, 14).Value = "`$/docena de atados"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 1167
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
